$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial for every data row (2-158).
# The whole column was bulk-updated from 2023-09-23 (45192) to 2023-10-03 (45202).
$ws.Range("C2:C158").Value = 45202
